$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "Play Arctic Magic for Free - Review and Game Information" "Play Arctic Magic Free: Winter-themed Slot with High Volatility"
Replace-Text "Stunning graphics and winter-themed visuals with a magical touch" "Top-notch graphics with a winter-themed and magical atmosphere"
Replace-Text "High volatility and a maximum jackpot of up to €150,000" "High volatility for a chance at a jackpot of up to 150,000 €"
Replace-Text "Wild and Scatter symbols add excitement to the gameplay" "Affordable betting range starting from 9 cents per spin"
Replace-Text "Free Spins feature with tripled winnings" "Lucrative Free Spins feature with tripled winnings"
Replace-Text "Low number of paylines compared to other slots" "Wins may be less frequent due to the high volatility"
Replace-Text "Less frequent wins due to the high volatility" "Limited number of paylines with only nine available"
Replace-Text "Read our review of Arctic Magic, a high volatility slot with stunning graphics and a winter theme. Play for free and take a shot at the jackpot of up to €150,000." "Discover the magic of Arctic Magic slot. Play for free and win big with high volatility."
